# Auto-generated Excel COM-interop script
# Applies cached numeric value updates (market price refresh) to the
# Hyperion_Profits workbook across all 8 crafting-class sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 10669.417
$ws.Range("J17").Value = 10669.417
$ws.Range("L17").Value = 32008.251
$ws.Range("N17").Value = -32344.251
# Row 132
$ws.Range("H132").Value = 45460084
$ws.Range("I132").Value = 47624800
$ws.Range("K132").Value = 142874400
$ws.Range("M132").Value = -142871870
# Row 138
$ws.Range("H138").Value = 3158.7334
$ws.Range("I138").Value = 1852.5294
$ws.Range("J138").Value = 3541.5862
$ws.Range("K138").Value = 5557.5882
$ws.Range("L138").Value = 10624.7586
$ws.Range("M138").Value = -417.5882000000001
$ws.Range("N138").Value = -20904.7586
# Row 141
$ws.Range("H141").Value = 14925.723
$ws.Range("I141").Value = 10830.2
$ws.Range("J141").Value = 35403.332
$ws.Range("K141").Value = 32490.6
$ws.Range("L141").Value = 106209.996
$ws.Range("M141").Value = -27310.6
$ws.Range("N141").Value = -116569.996

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 15144.083
$ws.Range("I32").Value = 9967.125
$ws.Range("J32").Value = 25498
$ws.Range("K32").Value = 9967.125
$ws.Range("L32").Value = 25498
$ws.Range("M32").Value = -9680.125
$ws.Range("N32").Value = -26072
# Row 45
$ws.Range("H45").Value = 5757499
$ws.Range("I45").Value = 9592047
$ws.Range("J45").Value = 5677.6
$ws.Range("K45").Value = 9592047
$ws.Range("L45").Value = 5677.6
$ws.Range("M45").Value = -9591670
$ws.Range("N45").Value = -6431.6
# Row 61
$ws.Range("H61").Value = 4108.8945
$ws.Range("I61").Value = 4004.7058
$ws.Range("K61").Value = 4004.7058
$ws.Range("M61").Value = -3792.7058
# Row 74
$ws.Range("H74").Value = 28991.176
$ws.Range("I74").Value = 1952.1538
$ws.Range("K74").Value = 1952.1538
$ws.Range("M74").Value = -1078.1538
# Row 77
$ws.Range("H77").Value = 28991.176
$ws.Range("I77").Value = 1952.1538
$ws.Range("K77").Value = 9760.769
$ws.Range("M77").Value = -5392.769
# Row 102
$ws.Range("H102").Value = 6414411.5
$ws.Range("I102").Value = 7578940.5
$ws.Range("K102").Value = 7578940.5
$ws.Range("M102").Value = -7577318.5
# Row 122
$ws.Range("H122").Value = 598775.4399999999
$ws.Range("I122").Value = 2307.3684
$ws.Range("J122").Value = 1307081.2
$ws.Range("K122").Value = 6922.1052
$ws.Range("L122").Value = 3921243.6
$ws.Range("M122").Value = -4472.1052
$ws.Range("N122").Value = -3926143.6
# Row 132
$ws.Range("H132").Value = 2598.7144
$ws.Range("I132").Value = 1736.2258
$ws.Range("J132").Value = 5029.364
$ws.Range("K132").Value = 5208.6774
$ws.Range("L132").Value = 15088.092
$ws.Range("M132").Value = -2678.6774
$ws.Range("N132").Value = -20148.092
# Row 136
$ws.Range("H136").Value = 4108.8945
$ws.Range("I136").Value = 4004.7058
$ws.Range("K136").Value = 12014.1174
$ws.Range("M136").Value = -9464.117400000001

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 6854106.5
# Row 105
$ws.Range("H105").Value = 4809387
$ws.Range("I105").Value = 6945558
$ws.Range("K105").Value = 6945558
$ws.Range("M105").Value = -6943811

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 732.125
$ws.Range("I22").Value = 451.16666
$ws.Range("J22").Value = 1575
$ws.Range("K22").Value = 451.16666
$ws.Range("L22").Value = 1575
$ws.Range("M22").Value = -101.16666
$ws.Range("N22").Value = -2275
# Row 31
$ws.Range("H31").Value = 21069.604
$ws.Range("I31").Value = 2083.0454
$ws.Range("K31").Value = 2083.0454
$ws.Range("M31").Value = -1788.0454
# Row 34
$ws.Range("H34").Value = 21069.604
$ws.Range("I34").Value = 2083.0454
$ws.Range("K34").Value = 2083.0454
$ws.Range("M34").Value = -1881.0454
# Row 58
$ws.Range("H58").Value = 6769.346
$ws.Range("I58").Value = 7795.9414
$ws.Range("J58").Value = 4830.222
$ws.Range("K58").Value = 7795.9414
$ws.Range("L58").Value = 4830.222
$ws.Range("M58").Value = -7592.9414
$ws.Range("N58").Value = -5236.222
# Row 86
$ws.Range("H86").Value = 16320.0625
$ws.Range("I86").Value = 11347
$ws.Range("K86").Value = 11347
$ws.Range("M86").Value = -10224
# Row 89
$ws.Range("H89").Value = 16320.0625
$ws.Range("I89").Value = 11347
$ws.Range("K89").Value = 56735
$ws.Range("M89").Value = -51119
# Row 97
$ws.Range("H97").Value = 32331.666
$ws.Range("J97").Value = 32331.666
$ws.Range("L97").Value = 32331.666
$ws.Range("N97").Value = -34313.666
# Row 116
$ws.Range("H116").Value = 42999.5
$ws.Range("J116").Value = 42999.5
$ws.Range("L116").Value = 42999.5
$ws.Range("N116").Value = -52177.5
# Row 122
$ws.Range("H122").Value = 2918.2307
$ws.Range("I122").Value = 2812.4546
$ws.Range("K122").Value = 8437.363799999999
$ws.Range("M122").Value = -5987.363799999999
# Row 132
$ws.Range("H132").Value = 46816.75
$ws.Range("I132").Value = 1658.8462
$ws.Range("K132").Value = 4976.5386
$ws.Range("M132").Value = -2446.5386
# Row 134
$ws.Range("H134").Value = 3774.2778
$ws.Range("I134").Value = 2787.182
$ws.Range("J134").Value = 5325.4287
$ws.Range("K134").Value = 8361.545999999998
$ws.Range("L134").Value = 15976.2861
$ws.Range("M134").Value = -5826.545999999998
$ws.Range("N134").Value = -21046.2861
# Row 136
$ws.Range("H136").Value = 6769.346
$ws.Range("I136").Value = 7795.9414
$ws.Range("J136").Value = 4830.222
$ws.Range("K136").Value = 23387.8242
$ws.Range("L136").Value = 14490.666
$ws.Range("M136").Value = -20837.8242
$ws.Range("N136").Value = -19590.666

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 467316.25
$ws.Range("J2").Value = 817296.5600000001
$ws.Range("L2").Value = 4903779.36
$ws.Range("N2").Value = -4904005.36
# Row 4
$ws.Range("H4").Value = 10003498
$ws.Range("I4").Value = 11330242
$ws.Range("J4").Value = 52917.5
$ws.Range("K4").Value = 33990726
$ws.Range("L4").Value = 158752.5
$ws.Range("M4").Value = -33990614
$ws.Range("N4").Value = -158976.5

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 68
$ws.Range("H68").Value = 2991
$ws.Range("I68").Value = 2991
$ws.Range("K68").Value = 2991
$ws.Range("M68").Value = -2180
# Row 71
$ws.Range("H71").Value = 2991
$ws.Range("I71").Value = 2991
$ws.Range("K71").Value = 8973
$ws.Range("M71").Value = -4917
# Row 102
$ws.Range("H102").Value = 8435689
$ws.Range("I102").Value = 27778130
$ws.Range("J102").Value = 1988208.2
$ws.Range("K102").Value = 27778130
$ws.Range("L102").Value = 1988208.2
$ws.Range("M102").Value = -27776508
$ws.Range("N102").Value = -1991452.2
# Row 107
$ws.Range("H107").Value = 643.75
$ws.Range("I107").Value = 762.75
$ws.Range("K107").Value = 762.75
$ws.Range("M107").Value = 1157.25
# Row 122
$ws.Range("H122").Value = 210012.19
$ws.Range("I122").Value = 264330.2
$ws.Range("J122").Value = 4810.8887
$ws.Range("K122").Value = 792990.6000000001
$ws.Range("L122").Value = 14432.6661
$ws.Range("M122").Value = -790540.6000000001
$ws.Range("N122").Value = -19332.6661
# Row 132
$ws.Range("H132").Value = 3192.2693
$ws.Range("I132").Value = 2909.5454
$ws.Range("J132").Value = 4747.25
$ws.Range("K132").Value = 8728.636200000001
$ws.Range("L132").Value = 14241.75
$ws.Range("M132").Value = -6198.636200000001
$ws.Range("N132").Value = -19301.75

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 9276.583000000001
$ws.Range("J7").Value = 11627
$ws.Range("L7").Value = 11627
$ws.Range("N7").Value = -11851
# Row 16
$ws.Range("H16").Value = 1551.0769
$ws.Range("I16").Value = 973.8889
$ws.Range("K16").Value = 973.8889
$ws.Range("M16").Value = -803.8889
# Row 40
$ws.Range("H40").Value = 5288.357
$ws.Range("I40").Value = 4239.4287
$ws.Range("K40").Value = 4239.4287
$ws.Range("M40").Value = -4103.4287
# Row 100
$ws.Range("H100").Value = 40607.37
$ws.Range("I100").Value = 3599.9167
$ws.Range("J100").Value = 70213.336
$ws.Range("K100").Value = 3599.9167
$ws.Range("L100").Value = 70213.336
$ws.Range("M100").Value = -3058.9167
$ws.Range("N100").Value = -71295.336
# Row 109
$ws.Range("H109").Value = 50285
$ws.Range("J109").Value = 50285
$ws.Range("L109").Value = 50285
$ws.Range("N109").Value = -53059
# Row 122
$ws.Range("H122").Value = 5434.4585
$ws.Range("I122").Value = 2874.889
$ws.Range("K122").Value = 8624.667000000001
$ws.Range("M122").Value = -6174.667000000001
# Row 126
$ws.Range("H126").Value = 9276.583000000001
$ws.Range("J126").Value = 11627
$ws.Range("L126").Value = 34881
$ws.Range("N126").Value = -39821
# Row 132
$ws.Range("H132").Value = 4125.3145
$ws.Range("I132").Value = 3988.3547
$ws.Range("K132").Value = 11965.0641
$ws.Range("M132").Value = -9435.0641
# Row 136
$ws.Range("H136").Value = 41155.605
$ws.Range("I136").Value = 146976.72
$ws.Range("J136").Value = 5881.905
$ws.Range("K136").Value = 440930.16
$ws.Range("L136").Value = 17645.715
$ws.Range("M136").Value = -438380.16
$ws.Range("N136").Value = -22745.715

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 117
$ws.Range("H117").Value = 27335.5
$ws.Range("I117").Value = 24671
$ws.Range("J117").Value = 30000
$ws.Range("K117").Value = 24671
$ws.Range("L117").Value = 30000
$ws.Range("M117").Value = -20082
$ws.Range("N117").Value = -39178
# Row 122
$ws.Range("H122").Value = 2291.1765
$ws.Range("I122").Value = 2099.0527
$ws.Range("J122").Value = 2534.5334
$ws.Range("K122").Value = 6297.158100000001
$ws.Range("L122").Value = 7603.600199999999
$ws.Range("M122").Value = -3847.158100000001
$ws.Range("N122").Value = -12503.6002

Write-Host "Applied all market value updates."